$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.836.09"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.301.07"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  -1.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.85"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.51"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "3.294.70"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.582"
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.58"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "633.10"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "3.831.68"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.11"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "65.871.18"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "3.293.57"
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.35"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.905"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.84"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.62"
$ws.Range("E24").Value = "  +8.06%  "
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.95"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.64"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.17"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.37"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.09"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "543.29"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").Value = "3.795.88"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.52"
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.74"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.24"
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.128"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.69"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -8.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.334"
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0416"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.59"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.99%  "
